$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'35.521.63"
$ws.Range('D3').Value = "'1.909.14"
$ws.Range('E3').Value = "'  +3.33%  "
$ws.Range('E4').Value = "'  +0.59%  "
$ws.Range('D5').Value = "'246.78"
$ws.Range('E5').Value = "'  +5.17%  "
$ws.Range('D6').Value = "'0.635"
$ws.Range('E6').Value = "'  +2.09%  "
$ws.Range('E7').Value = "'  +0.56%  "
$ws.Range('D8').Value = "'42.01"
$ws.Range('E8').Value = "'  -1.10%  "
$ws.Range('E9').Value = "'  +2.64%  "
$ws.Range('D10').Value = "'0.0706"
$ws.Range('E11').Value = "'  +1.25%  "
$ws.Range('E12').Value = "'  +3.32%  "
$ws.Range('D13').Value = "'12.40"
$ws.Range('E13').Value = "'  +9.23%  "
$ws.Range('D14').Value = "'1.905.15"
$ws.Range('E14').Value = "'  +2.94%  "
$ws.Range('E15').Value = "'  +2.40%  "
$ws.Range('D16').Value = "'4.87"
$ws.Range('E16').Value = "'  +3.84%  "
$ws.Range('D17').Value = "'35.544.39"
$ws.Range('E17').Value = "'  +1.53%  "
$ws.Range('D18').Value = "'72.07"
$ws.Range('E18').Value = "'  +2.98%  "
$ws.Range('D19').Value = "'0.0₃0822"
$ws.Range('E19').Value = "'  +3.73%  "
$ws.Range('D20').Value = "'243.64"
$ws.Range('E20').Value = "'  +1.21%  "
$ws.Range('D21').Value = "'12.56"
$ws.Range('E21').Value = "'  +3.63%  "
$ws.Range('E22').Value = "'  +2.44%  "
$ws.Range('E23').Value = "'  +0.51%  "
$ws.Range('E24').Value = "'  +0.84%  "
$ws.Range('D25').Value = "'172.59"
$ws.Range('E25').Value = "'  +1.06%  "
$ws.Range('E26').Value = "'  +20.89%  "
$ws.Range('D27').Value = "'8.52"
$ws.Range('E27').Value = "'  +8.32%  "
$ws.Range('D28').Value = "'18.04"
$ws.Range('E28').Value = "'  +2.28%  "
$ws.Range('E29').Value = "'  +0.98%  "
$ws.Range('E30').Value = "'  +25.87%  "
$ws.Range('E31').Value = "'  +2.91%  "
$ws.Range('D32').Value = "'4.11"
$ws.Range('E32').Value = "'  +3.31%  "
$ws.Range('B33').Value = 'BinanceUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D33').Value = "'1.01"
$ws.Range('E33').Value = "'  +0.55%  "
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = "'4.18"
$ws.Range('E34').Value = "'  +5.57%  "
$ws.Range('E35').Value = "'  +6.78%  "
$ws.Range('D36').Value = "'1.38"
$ws.Range('E36').Value = "'  +10.89%  "
$ws.Range('D37').Value = "'2.04"
$ws.Range('E37').Value = "'  +1.98%  "
$ws.Range('E38').Value = "'  +4.11%  "
$ws.Range('D39').Value = "'0.0205"
$ws.Range('E39').Value = "'  +2.10%  "
$ws.Range('D40').Value = "'91.29"
$ws.Range('E40').Value = "'  -0.12%  "
$ws.Range('D41').Value = "'15.79"
$ws.Range('E41').Value = "'  +5.59%  "
$ws.Range('D42').Value = "'1.354.45"
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = "'0.0607"
$ws.Range('E43').Value = "'  +13.96%  "
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').Value = "'49.57"
$ws.Range('E44').Value = "'  +43.41%  "
$ws.Range('D45').Value = "'2.37"
$ws.Range('E45').Value = "'  +1.80%  "
$ws.Range('D46').Value = "'12.67"
$ws.Range('E46').Value = "'  -1.95%  "
$ws.Range('E47').Value = "'  +1.33%  "
$ws.Range('E48').Value = "'  +0.48%  "
$ws.Range('E49').Value = "'  +4.46%  "
$ws.Range('D50').Value = "'2.097.09"
$ws.Range('E50').Value = "'  +3.38%  "
$ws.Range('E51').Value = "'  +2.16%  "
